$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-10-25 Saturday" "2025-10-26 Sunday"

Replace-Text "614÷6=102, 2" "153÷6=25, 3"
Replace-Text "596÷9=66, 2" "482÷9=53, 5"
Replace-Text "468÷2=234, 0" "536÷7=76, 4"
Replace-Text "255÷3=85, 0" "659÷6=109, 5"
Replace-Text "127÷2=63, 1" "913÷4=228, 1"
Replace-Text "591÷4=147, 3" "842÷2=421, 0"
Replace-Text "789÷4=197, 1" "660÷6=110, 0"
Replace-Text "950÷4=237, 2" "836÷4=209, 0"
Replace-Text "387÷9=43, 0" "229÷7=32, 5"
Replace-Text "511÷7=73, 0" "122÷9=13, 5"
Replace-Text "674÷4=168, 2" "304÷3=101, 1"
Replace-Text "180÷8=22, 4" "631÷2=315, 1"
Replace-Text "897÷7=128, 1" "479÷2=239, 1"
Replace-Text "225÷8=28, 1" "428÷2=214, 0"
Replace-Text "793÷8=99, 1" "976÷8=122, 0"
Replace-Text "281÷8=35, 1" "748÷2=374, 0"
Replace-Text "857÷3=285, 2" "579÷5=115, 4"
Replace-Text "868÷5=173, 3" "838÷4=209, 2"
Replace-Text "944÷3=314, 2" "735÷9=81, 6"
Replace-Text "949÷6=158, 1" "713÷4=178, 1"
Replace-Text "194÷7=27, 5" "393÷7=56, 1"
Replace-Text "257÷6=42, 5" "281÷6=46, 5"
Replace-Text "958÷6=159, 4" "860÷6=143, 2"
Replace-Text "259÷3=86, 1" "407÷5=81, 2"
Replace-Text "979÷7=139, 6" "429÷7=61, 2"
